$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '37.171.31'
$ws.Range('E2').Value = '  +1.58%  '
$ws.Range('D3').Value = '2.000.97'
$ws.Range('E3').Value = '  +2.09%  '
$ws.Range('E4').Value = '  +0.00%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '246.02'
$ws.Range('E5').Value = '  +0.61%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.626'
$ws.Range('E6').Value = '  +1.82%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '59.82'
$ws.Range('E7').Value = '  +1.94%  '
$ws.Range('E8').Value = '  -0.02%  '
$ws.Range('E9').Value = '  +2.59%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.0805'
$ws.Range('E10').Value = '  +1.79%  '
$ws.Range('E11').Value = '  +1.30%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '15.06'
$ws.Range('E12').Value = '  +6.01%  '
$ws.Range('D14').Value = '2.296.20'
$ws.Range('E14').Value = '  +2.10%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.844'
$ws.Range('E15').Value = '  +0.57%  '
$ws.Range('E16').Value = '  +2.58%  '
$ws.Range('D17').Value = '2.008.39'
$ws.Range('E17').Value = '  +2.62%  '
$ws.Range('D18').Value = '37.108.84'
$ws.Range('E18').Value = '  +1.56%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '70.25'
$ws.Range('E19').Value = '  +0.68%  '
$ws.Range('D20').Value = '0.0₃0864'
$ws.Range('E20').Value = '  +1.77%  '
$ws.Range('E21').Value = '  +2.22%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '230.24'
$ws.Range('E22').Value = '  +0.16%  '
$ws.Range('E23').Value = '  +0.10%  '
$ws.Range('E24').Value = '  +0.20%  '
$ws.Range('E25').Value = '  +0.57%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '9.42'
$ws.Range('E26').Value = '  +2.89%  '
$ws.Range('E27').Value = '  +1.69%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '163.48'
$ws.Range('E28').Value = '  +1.71%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '19.63'
$ws.Range('E29').Value = '  +1.00%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '1.34'
$ws.Range('E30').Value = '  +11.29%  '
$ws.Range('E31').Value = '  +0.91%  '
$ws.Range('E32').Value = '  +1.54%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.0654'
$ws.Range('E33').Value = '  +6.74%  '
$ws.Range('E34').Value = '  +1.85%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '2.40'
$ws.Range('E35').Value = '  +5.09%  '
$ws.Range('E36').Value = '  +0.07%  '
$ws.Range('E37').Value = '  +2.46%  '
$ws.Range('E38').Value = '  -6.15%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '5.36'
$ws.Range('E39').Value = '  -1.98%  '
$ws.Range('E40').Value = '  -0.12%  '
$ws.Range('E41').Value = '  +0.99%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.0214'
$ws.Range('E42').Value = '  +1.84%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '1.17'
$ws.Range('E43').Value = '  +0.47%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '16.67'
$ws.Range('E44').Value = '  +5.68%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '90.92'
$ws.Range('E45').Value = '  +3.18%  '
$ws.Range('D46').Value = '1.372.62'
$ws.Range('E46').Value = '  +0.06%  '
$ws.Range('B47').Value = 'FraxShare'
$ws.Range('C47').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '7.45'
$ws.Range('E47').Value = '  +4.51%  '
$ws.Range('B48').Value = 'ARBITRUM'
$ws.Range('C48').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '1.04'
$ws.Range('E48').Value = '  +1.92%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '2.03'
$ws.Range('E49').Value = '  +11.99%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '46.20'
$ws.Range('E51').Value = '  +4.94%  '
